# userTransactionReport.docx - add a "Failure Reason" column to the
# transaction table (header cell + ${failureReason} placeholder cell),
# and resize a couple of existing columns to make room for it.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# 1. Resize table / existing columns
# ---------------------------------------------------------------------
# Overall table width 15304 -> 16301 dxa (expressed in points for the OM).
$t.PreferredWidth = 16301 / 20.0

# Column 1 ("Tenant Name" / ${userName}): 988 -> 1135 dxa.
$t.Columns.Item(1).Width = 1135 / 20.0

# Column 8 ("Transaction Status" / ${transactionStatus}): 1276 -> 1129 dxa.
$t.Columns.Item(8).Width = 1129 / 20.0

# ---------------------------------------------------------------------
# 2. Insert the new "Failure Reason" column right after column 8
# ---------------------------------------------------------------------
$newCol = $t.Columns.Add($t.Columns.Item(9))
$newCol.Width = 997 / 20.0

$headerCell = $t.Cell(1, 9)
$headerCell.VerticalAlignment = 1   # wdCellAlignVerticalCenter

$dataCell = $t.Cell(2, 9)
$dataCell.VerticalAlignment = 1     # wdCellAlignVerticalCenter

$ooxmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        "<w:document $ooxmlNs>" +
        '<w:body>' + $bodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# -- Header row cell: "Failure Reason" -----------------------------------
$headerXml = New-PkgXml(
    '<w:p>' +
      '<w:pPr>' +
        '<w:jc w:val="center"/>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>' +
          '<w:b/>' +
          '<w:color w:val="000000"/>' +
          '<w:sz w:val="16"/>' +
          '<w:szCs w:val="16"/>' +
          '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '</w:rPr>' +
      '</w:pPr>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>' +
          '<w:b/>' +
          '<w:color w:val="000000"/>' +
          '<w:sz w:val="16"/>' +
          '<w:szCs w:val="16"/>' +
          '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '</w:rPr>' +
        '<w:t>Failure Reason</w:t>' +
      '</w:r>' +
    '</w:p>'
)
$headerCell.Range.Paragraphs.Item(1).Range.InsertXML($headerXml)

# -- Data row cell: ${failureReason} (with the relocated _GoBack mark) ---
$dataXml = New-PkgXml(
    '<w:p>' +
      '<w:pPr>' +
        '<w:jc w:val="center"/>' +
        '<w:rPr>' +
          '<w:sz w:val="12"/>' +
          '<w:szCs w:val="12"/>' +
        '</w:rPr>' +
      '</w:pPr>' +
      '<w:r>' +
        '<w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr>' +
        '<w:t>${</w:t>' +
      '</w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r>' +
        '<w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr>' +
        '<w:t>failureReason</w:t>' +
      '</w:r>' +
      '<w:bookmarkEnd w:id="0"/>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r>' +
        '<w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr>' +
        '<w:t>}</w:t>' +
      '</w:r>' +
    '</w:p>'
)
$dataCell.Range.Paragraphs.Item(1).Range.InsertXML($dataXml)

# ---------------------------------------------------------------------
# 3. The old _GoBack bookmark (sitting in its own empty paragraph right
#    after the table) moves into the new cell above, so clear it out of
#    its old home, leaving a plain empty paragraph behind.
# ---------------------------------------------------------------------
$tblEnd = $t.Range.End
$oldBookmarkPara = $d.Range($tblEnd + 2, $tblEnd + 3)
$emptyParaXml = New-PkgXml('<w:p/>')
$oldBookmarkPara.InsertXML($emptyParaXml)
